$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.199.43"
$ws.Range("E2").Value = "  -4.11%  "
$ws.Range("D3").Value = "2.240.23"
$ws.Range("E3").Value = "  -4.85%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "231.83"
$ws.Range("E5").Value = "  -3.71%  "
$ws.Range("D6").Value = "0.635"
$ws.Range("E6").Value = "  -5.89%  "
$ws.Range("D7").Value = "70.43"
$ws.Range("E7").Value = "  -4.27%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "0.560"
$ws.Range("E9").Value = "  -7.25%  "
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").Value = "58.32"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "35.70"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "6.85"
$ws.Range("E14").Value = "  -6.87%  "
$ws.Range("D15").Value = "2.578.13"
$ws.Range("E15").Value = "  -4.67%  "
$ws.Range("D16").Value = "15.06"
$ws.Range("E16").Value = "  -8.41%  "
$ws.Range("D17").Value = "0.866"
$ws.Range("E17").Value = "  -5.26%  "
$ws.Range("D18").Value = "2.241.42"
$ws.Range("E18").Value = "  -4.79%  "
$ws.Range("D19").Value = "42.096.55"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").Value = "73.53"
$ws.Range("E21").Value = "  -5.38%  "
$ws.Range("D22").Value = "6.21"
$ws.Range("E22").Value = "  -7.58%  "
$ws.Range("D23").Value = "238.62"
$ws.Range("E23").Value = "  -7.27%  "
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "3.64"
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("D27").Value = "2.35"
$ws.Range("E27").Value = "  -6.16%  "
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  -4.61%  "
$ws.Range("E29").Value = "  -5.05%  "
$ws.Range("D30").Value = "168.24"
$ws.Range("E30").Value = "  -5.50%  "
$ws.Range("D31").Value = "20.71"
$ws.Range("E31").Value = "  -8.83%  "
$ws.Range("E32").Value = "  -7.17%  "
$ws.Range("E33").Value = "  -7.07%  "
$ws.Range("D34").Value = "0.0720"
$ws.Range("E34").Value = "  -5.09%  "
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("D36").Value = "4.79"
$ws.Range("E36").Value = "  -8.07%  "
$ws.Range("D37").Value = "3.61"
$ws.Range("E37").Value = "  -5.74%  "
$ws.Range("D38").Value = "22.37"
$ws.Range("E38").Value = "  +17.02%  "
$ws.Range("D39").Value = "6.10"
$ws.Range("E39").Value = "  -4.48%  "
$ws.Range("E40").Value = "  -6.37%  "
$ws.Range("D41").Value = "0.0265"
$ws.Range("E41").Value = "  -5.01%  "
$ws.Range("D42").Value = "66.87"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "4.88"
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("D44").Value = "8.96"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("E45").Value = "  -9.77%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -7.22%  "
$ws.Range("D48").Value = "10.28"
$ws.Range("E48").Value = "  +7.89%  "
$ws.Range("D49").Value = "4.37"
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("D50").Value = "2.35"
$ws.Range("E50").Value = "  -6.35%  "
$ws.Range("E51").Value = "  -6.77%  "
